$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row values -----------------------------------------------------
# The sst entries for "Donor_ID" / "Code" / "Name" / "BusinessKey" were
# reshuffled, which (since the <c> cells still point at the same shared-
# string indices) changes what each header cell actually displays:
#   A2: Donor_ID  -> BusinessKey
#   B2: Code      -> Code        (unchanged)
#   C2: Name      -> Donor_ID
#   D2: BusinessKey -> Name
$ws.Range("A2").Value = "BusinessKey"
$ws.Range("B2").Value = "Code"
$ws.Range("C2").Value = "Donor_ID"
$ws.Range("D2").Value = "Name"

# --- Worksheet VBA CodeName -------------------------------------------------
# sheetPr/@codeName changed from "Sheet10" to "Sheet12" (the file was moved
# between template folders and Excel re-stamped the hidden code name).
try {
    $ws.CodeName = "Sheet12"
} catch {
    # CodeName is normally only editable from the VBA IDE Properties window;
    # keep going if this host doesn't expose a writable property for it.
}

# --- Workbook window size ---------------------------------------------------
# bookViews/workbookView@windowWidth/@windowHeight changed from 7470x2760 to
# 28800x12585 (the window was maximized/resized before the last save).
try {
    $win = $excel.ActiveWindow
    $win.Width = 28800
    $win.Height = 12585
} catch {
    try {
        $excel.Width = 28800
        $excel.Height = 12585
    } catch {
        # Window geometry isn't always settable from a headless host; ignore.
    }
}
